$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI TPM values for columns G:T (rows 2-26)
$data = @(
    @(14.734287,44.202861,0.4000023944294819,0.400002394429482,3,1,3.477616333333334,10.432849,0.4986219472823802,0.4986219472823803,51.240197131221,461.161774180989,0.199449972828043,0.1994499728280431),
    @(14.734287,44.202861,0.4000023944294819,0.400002394429482,3,1,1.578165333333333,4.734496,0.2262779433422874,0.2262779433422874,23.253140954784,209.278268593056,0.09051171914349362,0.09051171914349365),
    @(14.734287,44.202861,0.4000023944294819,0.400002394429482,2,0.6666666666666666,0.3701243333333333,1.110373,0.05306856712579453,0.05306856712579454,5.453518153017,49.081663377153,0.0212275539192595,0.02122755391925951),
    @(14.734287,44.202861,0.4000023944294819,0.400002394429482,3,1,0.775099,2.325297,0.1111339882471103,0.1111339882471103,11.420531119413,102.784780074717,0.04445386140134204,0.04445386140134205),
    @(14.734287,44.202861,0.4000023944294819,0.400002394429482,3,1,0.7734500000000001,2.32035,0.1108975540024274,0.1108975540024274,11.39623428015,102.56610852135,0.04435928713734375,0.04435928713734377),
    @(14.452944,43.358832,0.3923645715978801,0.3923645715978802,3,1,3.477616333333334,10.432849,0.4986219472823802,0.4986219472823803,50.26179411915201,452.3561470723681,0.1956415867347519,0.195641586734752),
    @(14.452944,43.358832,0.3923645715978801,0.3923645715978802,3,1,1.578165333333333,4.734496,0.2262779433422874,0.2262779433422874,22.809135185408,205.282216668672,0.088783448301546,0.08878344830154603),
    @(14.452944,43.358832,0.3923645715978801,0.3923645715978802,2,0.6666666666666666,0.3701243333333333,1.110373,0.05306856712579453,0.05306856712579454,5.349386262704,48.144476364336,0.02082222560562572,0.02082222560562573),
    @(14.452944,43.358832,0.3923645715978801,0.3923645715978802,3,1,0.775099,2.325297,0.1111339882471103,0.1111339882471103,11.202462441456,100.822161973104,0.04360503968854129,0.04360503968854131),
    @(14.452944,43.358832,0.3923645715978801,0.3923645715978802,3,1,0.7734500000000001,2.32035,0.1108975540024274,0.1108975540024274,11.1786295368,100.6076658312,0.04351227126741521,0.04351227126741523),
    @(1.259379333333333,3.778138,0.03418928576783783,0.03418928576783784,3,1,3.477616333333334,10.432849,0.4986219472823802,0.4986219472823803,4.379638139462445,39.416743255162,0.01704752824575307,0.01704752824575308),
    @(1.259379333333333,3.778138,0.03418928576783783,0.03418928576783784,3,1,1.578165333333333,4.734496,0.2262779433422874,0.2262779433422874,1.987508805383111,17.887579248448,0.007736281267888083,0.007736281267888085),
    @(1.259379333333333,3.778138,0.03418928576783783,0.03418928576783784,2,0.6666666666666666,0.3701243333333333,1.110373,0.05306856712579453,0.05306856712579454,0.4661269361637778,4.195142425474001,0.001814376406753474,0.001814376406753474),
    @(1.259379333333333,3.778138,0.03418928576783783,0.03418928576783784,3,1,0.775099,2.325297,0.1111339882471103,0.1111339882471103,0.9761436618873334,8.785292956986,0.003799591682699987,0.003799591682699988),
    @(1.259379333333333,3.778138,0.03418928576783783,0.03418928576783784,3,1,0.7734500000000001,2.32035,0.1108975540024274,0.1108975540024274,0.9740669453666668,8.766602508300002,0.00379150816474322,0.003791508164743221),
    @(4.524255666666667,13.572767,0.1228232556945456,0.1228232556945456,3,1,3.477616333333334,10.432849,0.4986219472823802,0.4986219472823803,15.73362540257589,141.602628623183,0.06124237092597602,0.06124237092597604),
    @(4.524255666666667,13.572767,0.1228232556945456,0.1228232556945456,3,1,1.578165333333333,4.734496,0.2262779433422874,0.2262779433422874,7.140023452270223,64.26021107043201,0.02779219369316567,0.02779219369316567),
    @(4.524255666666667,13.572767,0.1228232556945456,0.1228232556945456,2,0.6666666666666666,0.3701243333333333,1.110373,0.05306856712579453,0.05306856712579454,1.674537112454556,15.070834012091,0.006518054189434618,0.00651805418943462),
    @(4.524255666666667,13.572767,0.1228232556945456,0.1228232556945456,3,1,0.775099,2.325297,0.1111339882471103,0.1111339882471103,3.506746042977667,31.560714386799,0.01364983825482946,0.01364983825482946),
    @(4.524255666666667,13.572767,0.1228232556945456,0.1228232556945456,3,1,0.7734500000000001,2.32035,0.1108975540024274,0.1108975540024274,3.499285545383334,31.49356990845001,0.01362079863113982,0.01362079863113983),
    @(1.864631,5.593893,0.05062049251025444,0.05062049251025445,3,1,3.477616333333334,10.432849,0.4986219472823802,0.4986219472823803,6.484471221239667,58.360240991157,0.02524048854785621,0.02524048854785622),
    @(1.864631,5.593893,0.05062049251025444,0.05062049251025445,3,1,1.578165333333333,4.734496,0.2262779433422874,0.2262779433422874,2.942696003658667,26.484264032928,0.01145430093619404,0.01145430093619404),
    @(1.864631,5.593893,0.05062049251025444,0.05062049251025445,2,0.6666666666666666,0.3701243333333333,1.110373,0.05306856712579453,0.05306856712579454,0.6901453057876666,6.211307752089,0.002686357004721217,0.002686357004721218),
    @(1.864631,5.593893,0.05062049251025444,0.05062049251025445,3,1,0.775099,2.325297,0.1111339882471103,0.1111339882471103,1.445273623469,13.007462611221,0.005625657219697554,0.005625657219697555),
    @(1.864631,5.593893,0.05062049251025444,0.05062049251025445,3,1,0.7734500000000001,2.32035,0.1108975540024274,0.1108975540024274,1.44219884695,12.97978962255,0.005613688801785414,0.005613688801785416)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($i + 2, $j + 7).Value = $row[$j]
    }
}

Write-Output "done"